$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching the style/formatting of the other header
# cells (copy the format from G1, the last existing header, then set the text)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the value 0 in H2 for the new "Save" column
$ws.Range("H2").Value = 0
